# Translation patch v2: the former "column D" (and the occasional
# "column C") English/localized text now lives in column B, directly next
# to the Japanese source text in column A. Columns C and D are retired.
#
# For every row 1..71:
#   - new B = old D (if present), else old C (if present), else a copy of A
#   - old C and D are cleared
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$maxRow = 71
for ($r = 1; $r -le $maxRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    $d = $ws.Cells.Item($r, 4).Value()

    if ($d -ne $null -and $d -ne "") {
        $newB = $d
    } elseif ($c -ne $null -and $c -ne "") {
        $newB = $c
    } else {
        $newB = $a
    }

    $ws.Cells.Item($r, 2).Value = $newB
    $ws.Cells.Item($r, 3).Value = $null
    $ws.Cells.Item($r, 4).Value = $null
}
